$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1322.2222
$ws.Range("I12").Value = 1614
$ws.Range("K12").Value = 1614
$ws.Range("M12").Value = -1444

$ws.Range("H62").Value = 2149.6667
$ws.Range("I62").Value = 2149.6667
$ws.Range("K62").Value = 2149.6667
$ws.Range("M62").Value = -1525.6667

$ws.Range("H65").Value = 2149.6667
$ws.Range("I65").Value = 2149.6667
$ws.Range("K65").Value = 10748.3335
$ws.Range("M65").Value = -7628.333500000001

$ws.Range("H86").Value = 21865.2
$ws.Range("I86").Value = 1920
$ws.Range("K86").Value = 1920
$ws.Range("M86").Value = -797

$ws.Range("H89").Value = 21865.2
$ws.Range("I89").Value = 1920
$ws.Range("K89").Value = 9600
$ws.Range("M89").Value = -3984

$ws.Range("H113").Value = 25004574
$ws.Range("I113").Value = 142859600
$ws.Range("J113").Value = 5021.0303
$ws.Range("K113").Value = 142859600
$ws.Range("L113").Value = 5021.0303
$ws.Range("M113").Value = -142856346
$ws.Range("N113").Value = -11529.0303

$ws.Range("H135").Value = 25008058
$ws.Range("I135").Value = 1071.625
$ws.Range("J135").Value = 125036000
$ws.Range("K135").Value = 9644.625
$ws.Range("L135").Value = 1125324000
$ws.Range("M135").Value = -7109.625
$ws.Range("N135").Value = -1125329070

$ws.Range("H138").Value = 2559.7222
$ws.Range("J138").Value = 2873.4211
$ws.Range("L138").Value = 8620.263300000001
$ws.Range("N138").Value = -18900.2633

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").Value = 0

$ws.Range("H44").Value = 29666.666
$ws.Range("J44").Value = 29666.666
$ws.Range("L44").Value = 29666.666
$ws.Range("N44").Value = -30642.666

$ws.Range("H55").Value = 26993.334
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").ClearContents()

$ws.Range("H61").Value = 431658.6
$ws.Range("I61").Value = 785105.25
$ws.Range("J61").Value = 3802.1052
$ws.Range("K61").Value = 785105.25
$ws.Range("L61").Value = 3802.1052
$ws.Range("M61").Value = -784893.25
$ws.Range("N61").Value = -4226.1052

$ws.Range("H63").Value = 3474256
$ws.Range("I63").Value = 2288.125
$ws.Range("J63").Value = 31250000
$ws.Range("K63").Value = 2288.125
$ws.Range("L63").Value = 31250000
$ws.Range("M63").Value = -1602.125
$ws.Range("N63").Value = -31251372

$ws.Range("H66").Value = 3474256
$ws.Range("I66").Value = 2288.125
$ws.Range("J66").Value = 31250000
$ws.Range("K66").Value = 11440.625
$ws.Range("L66").Value = 156250000
$ws.Range("M66").Value = -8008.625
$ws.Range("N66").Value = -156256864

$ws.Range("H80").Value = 41336
$ws.Range("J80").Value = 41336
$ws.Range("L80").Value = 41336
$ws.Range("N80").Value = -43332

$ws.Range("H83").Value = 41336
$ws.Range("J83").Value = 41336
$ws.Range("L83").Value = 124008
$ws.Range("N83").Value = -133992

$ws.Range("H122").Value = 2235.2888
$ws.Range("I122").Value = 2134.543
$ws.Range("K122").Value = 6403.629000000001
$ws.Range("M122").Value = -3953.629000000001

$ws.Range("H132").Value = 15128
$ws.Range("I132").Value = 1468.9474
$ws.Range("J132").Value = 29545.889
$ws.Range("K132").Value = 4406.8422
$ws.Range("L132").Value = 88637.667
$ws.Range("M132").Value = -1876.8422
$ws.Range("N132").Value = -93697.667

$ws.Range("H136").Value = 431658.6
$ws.Range("I136").Value = 785105.25
$ws.Range("J136").Value = 3802.1052
$ws.Range("K136").Value = 2355315.75
$ws.Range("L136").Value = 11406.3156
$ws.Range("M136").Value = -2352765.75
$ws.Range("N136").Value = -16506.3156

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 19701.445
$ws.Range("I58").Value = 1208.1428
$ws.Range("J58").Value = 84428
$ws.Range("K58").Value = 1208.1428
$ws.Range("L58").Value = 84428
$ws.Range("M58").Value = -1005.1428
$ws.Range("N58").Value = -84834

$ws.Range("H86").Value = 12750.277
$ws.Range("I86").Value = 6272.4546
$ws.Range("J86").Value = 22929.715
$ws.Range("K86").Value = 6272.4546
$ws.Range("L86").Value = 22929.715
$ws.Range("M86").Value = -5149.4546
$ws.Range("N86").Value = -25175.715

$ws.Range("H89").Value = 12750.277
$ws.Range("I89").Value = 6272.4546
$ws.Range("J89").Value = 22929.715
$ws.Range("K89").Value = 31362.273
$ws.Range("L89").Value = 114648.575
$ws.Range("M89").Value = -25746.273
$ws.Range("N89").Value = -125880.575

$ws.Range("H99").Value = 5576.143
$ws.Range("I99").Value = 4160
$ws.Range("J99").Value = 6863.5454
$ws.Range("K99").Value = 4160
$ws.Range("L99").Value = 6863.5454
$ws.Range("M99").Value = -2662
$ws.Range("N99").Value = -9859.545399999999

$ws.Range("H107").Value = 337.35294
$ws.Range("I107").Value = 363.27274
$ws.Range("J107").Value = 289.83334
$ws.Range("K107").Value = 363.27274
$ws.Range("L107").Value = 289.83334
$ws.Range("M107").Value = 1556.72726
$ws.Range("N107").Value = -4129.83334

$ws.Range("H126").Value = 5576.143
$ws.Range("I126").Value = 4160
$ws.Range("J126").Value = 6863.5454
$ws.Range("K126").Value = 12480
$ws.Range("L126").Value = 20590.6362
$ws.Range("M126").Value = -10010
$ws.Range("N126").Value = -25530.6362

$ws.Range("H134").Value = 3242.9744
$ws.Range("I134").Value = 574.3333
$ws.Range("J134").Value = 35266.668
$ws.Range("K134").Value = 1722.9999
$ws.Range("L134").Value = 105800.004
$ws.Range("M134").Value = 812.0001
$ws.Range("N134").Value = -110870.004

$ws.Range("H136").Value = 19701.445
$ws.Range("I136").Value = 1208.1428
$ws.Range("J136").Value = 84428
$ws.Range("K136").Value = 3624.4284
$ws.Range("L136").Value = 253284
$ws.Range("M136").Value = -1074.4284
$ws.Range("N136").Value = -258384

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 246.83333
$ws.Range("I33").Value = 145.5
$ws.Range("J33").Value = 297.5
$ws.Range("K33").Value = 873
$ws.Range("L33").Value = 1785
$ws.Range("M33").Value = -590
$ws.Range("N33").Value = -2351

$ws.Range("H107").Value = 5061.5415
$ws.Range("J107").Value = 998.85
$ws.Range("L107").Value = 2996.55
$ws.Range("N107").Value = -6836.55

$ws.Range("H109").Value = 1840.875
$ws.Range("I109").Value = 806.75
$ws.Range("J109").Value = 2875
$ws.Range("K109").Value = 2420.25
$ws.Range("L109").Value = 8625
$ws.Range("M109").Value = -1380.25
$ws.Range("N109").Value = -10705

$ws.Range("H121").Value = 4546.7407
$ws.Range("I121").Value = 676.6667
$ws.Range("J121").Value = 5030.5
$ws.Range("K121").Value = 2030.0001
$ws.Range("L121").Value = 15091.5
$ws.Range("M121").Value = -720.0001
$ws.Range("N121").Value = -17711.5

$ws.Range("H131").Value = 824.51
$ws.Range("J131").Value = 832.1531
$ws.Range("L131").Value = 2496.4593
$ws.Range("N131").Value = -12576.4593

$ws.Range("H132").Value = 800.4
$ws.Range("J132").Value = 801.25
$ws.Range("L132").Value = 7211.25
$ws.Range("N132").Value = -12271.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8932367
$ws.Range("J70").Value = 10420428
$ws.Range("L70").Value = 10420428
$ws.Range("N70").Value = -10420968

$ws.Range("H73").Value = 8932367
$ws.Range("J73").Value = 10420428
$ws.Range("L73").Value = 10420428
$ws.Range("N73").Value = -10422300

$ws.Range("H102").Value = 1813.05
$ws.Range("I102").Value = 1746.1875
$ws.Range("J102").Value = 2080.5
$ws.Range("K102").Value = 1746.1875
$ws.Range("L102").Value = 2080.5
$ws.Range("M102").Value = -124.1875
$ws.Range("N102").Value = -5324.5

$ws.Range("H126").Value = 5873.684
$ws.Range("I126").Value = 4410
$ws.Range("J126").Value = 7500
$ws.Range("K126").Value = 13230
$ws.Range("L126").Value = 22500
$ws.Range("M126").Value = -10760
$ws.Range("N126").Value = -27440

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5206.8696
$ws.Range("I7").Value = 5093.5454
$ws.Range("K7").Value = 5093.5454
$ws.Range("M7").Value = -4981.5454

$ws.Range("H40").Value = 105015.91
$ws.Range("I40").Value = 141774.38
$ws.Range("K40").Value = 141774.38
$ws.Range("M40").Value = -141638.38

$ws.Range("H68").Value = 3521.2307
$ws.Range("I68").Value = 2189.111
$ws.Range("J68").Value = 4226.4707
$ws.Range("K68").Value = 2189.111
$ws.Range("L68").Value = 4226.4707
$ws.Range("M68").Value = -1440.111
$ws.Range("N68").Value = -5724.4707

$ws.Range("H71").Value = 3521.2307
$ws.Range("I71").Value = 2189.111
$ws.Range("J71").Value = 4226.4707
$ws.Range("K71").Value = 10945.555
$ws.Range("L71").Value = 21132.3535
$ws.Range("M71").Value = -7201.555
$ws.Range("N71").Value = -28620.3535

$ws.Range("H93").Value = 1779.3158
$ws.Range("I93").Value = 1876.6875
$ws.Range("K93").Value = 1876.6875
$ws.Range("M93").Value = -628.6875

$ws.Range("H126").Value = 5206.8696
$ws.Range("I126").Value = 5093.5454
$ws.Range("K126").Value = 15280.6362
$ws.Range("M126").Value = -12810.6362

$ws.Range("H132").Value = 1650.4865
$ws.Range("I132").Value = 1135.3334
$ws.Range("K132").Value = 3406.0002
$ws.Range("M132").Value = -876.0001999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1160.6333
$ws.Range("I136").Value = 691.8182
$ws.Range("J136").Value = 2449.875
$ws.Range("K136").Value = 2075.4546
$ws.Range("L136").Value = 7349.625
$ws.Range("M136").Value = 474.5454
$ws.Range("N136").Value = -12449.625
